$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row (row 1) ---
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"

# New header cells H1:M1 need the bold/centered/bordered "header" style
# to match the existing B1:G1 formatting (style index 1 in the sheet).
$newHeader = $ws.Range("H1:M1")
$newHeader.Font.Bold = $true
$newHeader.HorizontalAlignment = -4108
$newHeader.VerticalAlignment = -4160
$newHeader.Borders.LineStyle = 1
$newHeader.Borders.Weight = 2

$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows 2-19 ---
# Row 2 (property index 48)
$ws.Range("B2").Value = "臺灣中小企業銀行安平分行"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "許添財"
$ws.Range("F2").Value = 2165459
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "許添財"
$ws.Range("K2").Value = 639
$ws.Range("L2").Value = "tmpb8d31"
$ws.Range("M2").Value = 48

# Row 3 (property index 49)
$ws.Range("B3").Value = "臺灣銀行群賢分行"
$ws.Range("C3").Value = "綜合存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "洪淑頁"
$ws.Range("F3").Value = 2253656
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("J3").Value = "許添財"
$ws.Range("K3").Value = 639
$ws.Range("L3").Value = "tmpb8d31"
$ws.Range("M3").Value = 49

# Row 4 (property index 50)
$ws.Range("B4").Value = "臺灣銀行南都分行"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "許添財"
$ws.Range("F4").Value = 1779236
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("J4").Value = "許添財"
$ws.Range("K4").Value = 639
$ws.Range("L4").Value = "tmpb8d31"
$ws.Range("M4").Value = 50

# Row 5 (property index 51)
$ws.Range("B5").Value = "臺灣銀行南都分行"
$ws.Range("C5").Value = "定期存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "許添財"
$ws.Range("F5").Value = 3000000
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("J5").Value = "許添財"
$ws.Range("K5").Value = 639
$ws.Range("L5").Value = "tmpb8d31"
$ws.Range("M5").Value = 51

# Row 6 (property index 52)
$ws.Range("B6").Value = "臺灣銀行南都分行"
$ws.Range("C6").Value = "綜合存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "洪淑貞"
$ws.Range("F6").Value = 1000
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("J6").Value = "許添財"
$ws.Range("K6").Value = 639
$ws.Range("L6").Value = "tmpb8d31"
$ws.Range("M6").Value = 52

# Row 7 (property index 53)
$ws.Range("B7").Value = "台北富邦商業銀行駐立分行"
$ws.Range("C7").Value = "活期儲蓄存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "許添財"
$ws.Range("F7").Value = 254
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("J7").Value = "許添財"
$ws.Range("K7").Value = 639
$ws.Range("L7").Value = "tmpb8d31"
$ws.Range("M7").Value = 53

# Row 8 (property index 54)
$ws.Range("B8").Value = "華南商業銀行台南分行"
$ws.Range("C8").Value = "活期存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "許添財"
$ws.Range("F8").Value = 1210
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("J8").Value = "許添財"
$ws.Range("K8").Value = 639
$ws.Range("L8").Value = "tmpb8d31"
$ws.Range("M8").Value = 54

# Row 9 (property index 55)
$ws.Range("B9").Value = "華南商業銀行東台南分行"
$ws.Range("C9").Value = "活期存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "許添財"
$ws.Range("F9").Value = 1845
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("J9").Value = "許添財"
$ws.Range("K9").Value = 639
$ws.Range("L9").Value = "tmpb8d31"
$ws.Range("M9").Value = 55

# Row 10 (property index 56)
$ws.Range("B10").Value = "第一商業銀行運河分行"
$ws.Range("C10").Value = "活期儲蓄存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "許添財"
$ws.Range("F10").Value = 11567
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("J10").Value = "許添財"
$ws.Range("K10").Value = 639
$ws.Range("L10").Value = "tmpb8d31"
$ws.Range("M10").Value = 56

# Row 11 (property index 57)
$ws.Range("B11").Value = "合作金庫商業銀行南興分行"
$ws.Range("C11").Value = "支票存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "許添財"
$ws.Range("F11").Value = 1977
$ws.Range("G11").Value = "deposit"
$ws.Range("H11").Value = "normal"
$ws.Range("J11").Value = "許添財"
$ws.Range("K11").Value = 639
$ws.Range("L11").Value = "tmpb8d31"
$ws.Range("M11").Value = 57

# Row 12 (property index 58)
$ws.Range("B12").Value = "金城商業銀行府城分行"
$ws.Range("C12").Value = "活期存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "許添財"
$ws.Range("F12").Value = 388
$ws.Range("G12").Value = "deposit"
$ws.Range("H12").Value = "normal"
$ws.Range("J12").Value = "許添財"
$ws.Range("K12").Value = 639
$ws.Range("L12").Value = "tmpb8d31"
$ws.Range("M12").Value = 58

# Row 13 (property index 59)
$ws.Range("B13").Value = "臺灣銀行(註2)"
$ws.Range("C13").Value = "綜合存款"
$ws.Range("D13").Value = "美金"
$ws.Range("E13").Value = "許添財"
$ws.Range("F13").Value = 10991
$ws.Range("G13").Value = "deposit"
$ws.Range("H13").Value = "normal"
$ws.Range("J13").Value = "許添財"
$ws.Range("K13").Value = 639
$ws.Range("L13").Value = "tmpb8d31"
$ws.Range("M13").Value = 59

# Row 14 (property index 60)
$ws.Range("B14").Value = "美商花旗銀行"
$ws.Range("C14").Value = "支票存款"
$ws.Range("D14").Value = "美金"
$ws.Range("E14").Value = "許添財洪淑貞"
$ws.Range("F14").Value = 59475
$ws.Range("G14").Value = "deposit"
$ws.Range("H14").Value = "normal"
$ws.Range("J14").Value = "許添財"
$ws.Range("K14").Value = 639
$ws.Range("L14").Value = "tmpb8d31"
$ws.Range("M14").Value = 60

# Row 15 (property index 61)
$ws.Range("B15").Value = "日商三菱東京曰聯銀行"
$ws.Range("C15").Value = "其他存款"
$ws.Range("D15").Value = "美金"
$ws.Range("E15").Value = "洪淑貞"
$ws.Range("F15").Value = 3001588
$ws.Range("G15").Value = "deposit"
$ws.Range("H15").Value = "normal"
$ws.Range("J15").Value = "許添財"
$ws.Range("K15").Value = 639
$ws.Range("L15").Value = "tmpb8d31"
$ws.Range("M15").Value = 61

# Row 16 (property index 62)
$ws.Range("B16").Value = "美商花旗銀行"
$ws.Range("C16").Value = "其他存款"
$ws.Range("D16").Value = "美金"
$ws.Range("E16").Value = "許添財"
$ws.Range("F16").Value = 492278
$ws.Range("G16").Value = "deposit"
$ws.Range("H16").Value = "normal"
$ws.Range("J16").Value = "許添財"
$ws.Range("K16").Value = 639
$ws.Range("L16").Value = "tmpb8d31"
$ws.Range("M16").Value = 62

# Row 17 (property index 63)
$ws.Range("B17").Value = "美商花旗銀行"
$ws.Range("C17").Value = "其他存款"
$ws.Range("D17").Value = "美金"
$ws.Range("E17").Value = "洪淑貞"
$ws.Range("F17").Value = 145381
$ws.Range("G17").Value = "deposit"
$ws.Range("H17").Value = "normal"
$ws.Range("J17").Value = "許添財"
$ws.Range("K17").Value = 639
$ws.Range("L17").Value = "tmpb8d31"
$ws.Range("M17").Value = 63

# Row 18 (property index 64)
$ws.Range("B18").Value = "美商摩根大通銀行"
$ws.Range("C18").Value = "其他存款"
$ws.Range("D18").Value = "美金"
$ws.Range("E18").Value = "許添財"
$ws.Range("F18").Value = 124909
$ws.Range("G18").Value = "deposit"
$ws.Range("H18").Value = "normal"
$ws.Range("J18").Value = "許添財"
$ws.Range("K18").Value = 639
$ws.Range("L18").Value = "tmpb8d31"
$ws.Range("M18").Value = 64

# Row 19 (property index 65)
$ws.Range("B19").Value = "美商摩根大通銀行"
$ws.Range("C19").Value = "其他存款"
$ws.Range("D19").Value = "美金"
$ws.Range("E19").Value = "洪淑貞"
$ws.Range("F19").Value = 124909
$ws.Range("G19").Value = "deposit"
$ws.Range("H19").Value = "normal"
$ws.Range("J19").Value = "許添財"
$ws.Range("K19").Value = 639
$ws.Range("L19").Value = "tmpb8d31"
$ws.Range("M19").Value = 65

# --- Date column (I2:I19): force text storage, matching the shared-string
#     "date" values used elsewhere in the workbook, instead of letting Excel
#     auto-convert the "yyyy-mm-dd" text into a real date serial number. ---
$dateRange = $ws.Range("I2:I19")
$dateRange.NumberFormat = "@"
$ws.Range("I2").Value = "2012-03-22"
$ws.Range("I3").Value = "2012-03-22"
$ws.Range("I4").Value = "2012-03-22"
$ws.Range("I5").Value = "2012-03-22"
$ws.Range("I6").Value = "2012-03-22"
$ws.Range("I7").Value = "2012-03-22"
$ws.Range("I8").Value = "2012-03-22"
$ws.Range("I9").Value = "2012-03-22"
$ws.Range("I10").Value = "2012-03-22"
$ws.Range("I11").Value = "2012-03-22"
$ws.Range("I12").Value = "2012-03-22"
$ws.Range("I13").Value = "2012-03-22"
$ws.Range("I14").Value = "2012-03-22"
$ws.Range("I15").Value = "2012-03-22"
$ws.Range("I16").Value = "2012-03-22"
$ws.Range("I17").Value = "2012-03-22"
$ws.Range("I18").Value = "2012-03-22"
$ws.Range("I19").Value = "2012-03-22"
$dateRange.Style = "Normal"

